$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Experimental = false (force literal text "false", not boolean FALSE).
# A bare Value = "false" gets auto-typed to a Boolean by Excel, so enter it
# with a leading apostrophe to force text, then restore the original
# (non-quote-prefixed) cell format by pasting formats from a sibling cell.
$ws.Range("B7").Value = "'false"
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# Date updated
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Description added
$ws.Range("B17").Value = "Current training status categories based on fitness trends"
